$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 79.5
$ws.Range("I9").Value = 79.5
$ws.Range("K9").Value = 79.5
$ws.Range("M9").Value = 89.5
$ws.Range("H38").Value = 2509.9
$ws.Range("I38").Value = 637.75
$ws.Range("J38").Value = 9998.5
$ws.Range("K38").Value = 1913.25
$ws.Range("L38").Value = 29995.5
$ws.Range("M38").Value = -1541.25
$ws.Range("N38").Value = -30739.5
$ws.Range("H39").Value = 360.4
$ws.Range("I39").Value = 169.25
$ws.Range("J39").Value = 1125
$ws.Range("K39").Value = 507.75
$ws.Range("L39").Value = 3375
$ws.Range("M39").Value = -211.75
$ws.Range("N39").Value = -3967
$ws.Range("H41").Value = 804.375
$ws.Range("I41").Value = 612.3333
$ws.Range("J41").Value = 919.6
$ws.Range("K41").Value = 612.3333
$ws.Range("L41").Value = 919.6
$ws.Range("M41").Value = -172.3333
$ws.Range("N41").Value = -1799.6
$ws.Range("I64").Value = 3163.3333
$ws.Range("J64").Value = 3250
$ws.Range("K64").Value = 3163.3333
$ws.Range("L64").Value = 3250
$ws.Range("M64").Value = -2915.3333
$ws.Range("N64").Value = -3746
$ws.Range("I67").Value = 3163.3333
$ws.Range("J67").Value = 3250
$ws.Range("K67").Value = 3163.3333
$ws.Range("L67").Value = 3250
$ws.Range("M67").Value = -2305.3333
$ws.Range("N67").Value = -4966
$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 5000
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = -14126
$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 5000
$ws.Range("K72").Value = 45000
$ws.Range("M72").Value = -40632
$ws.Range("H107").Value = 451.4
$ws.Range("I107").Value = 451.4
$ws.Range("K107").Value = 451.4
$ws.Range("M107").Value = 1468.6
$ws.Range("H138").Value = 4256.067
$ws.Range("I138").Value = 1014.5
$ws.Range("K138").Value = 3043.5
$ws.Range("M138").Value = 2096.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4059
$ws.Range("I2").Value = 3638.6667
$ws.Range("J2").Value = 4899.6665
$ws.Range("K2").Value = 3638.6667
$ws.Range("L2").Value = 4899.6665
$ws.Range("M2").Value = -3525.6667
$ws.Range("N2").Value = -5125.6665
$ws.Range("H25").Value = 250
$ws.Range("I25").Value = 250
$ws.Range("K25").Value = 250
$ws.Range("M25").Value = 152
$ws.Range("H26").Value = 2485.7856
$ws.Range("I26").Value = 975.125
$ws.Range("K26").Value = 975.125
$ws.Range("M26").Value = -645.125
$ws.Range("H32").Value = 5040.551
$ws.Range("I32").Value = 3449.6943
$ws.Range("K32").Value = 3449.6943
$ws.Range("M32").Value = -3162.6943
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H61").Value = 1976.8182
$ws.Range("I61").Value = 1773.0555
$ws.Range("J61").Value = 2893.75
$ws.Range("K61").Value = 1773.0555
$ws.Range("L61").Value = 2893.75
$ws.Range("M61").Value = -1561.0555
$ws.Range("N61").Value = -3317.75
$ws.Range("H116").Value = 4059
$ws.Range("I116").Value = 3638.6667
$ws.Range("J116").Value = 4899.6665
$ws.Range("K116").Value = 3638.6667
$ws.Range("L116").Value = 4899.6665
$ws.Range("M116").Value = -1344.6667
$ws.Range("N116").Value = -9487.666499999999
$ws.Range("H132").Value = 2069.5293
$ws.Range("I132").Value = 2072.7273
$ws.Range("J132").Value = 2063.6667
$ws.Range("K132").Value = 6218.1819
$ws.Range("L132").Value = 6191.000100000001
$ws.Range("M132").Value = -3688.1819
$ws.Range("N132").Value = -11251.0001
$ws.Range("H136").Value = 1976.8182
$ws.Range("I136").Value = 1773.0555
$ws.Range("J136").Value = 2893.75
$ws.Range("K136").Value = 5319.166499999999
$ws.Range("L136").Value = 8681.25
$ws.Range("M136").Value = -2769.166499999999
$ws.Range("N136").Value = -13781.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4059
$ws.Range("I3").Value = 3638.6667
$ws.Range("J3").Value = 4899.6665
$ws.Range("K3").Value = 3638.6667
$ws.Range("L3").Value = 4899.6665
$ws.Range("M3").Value = -3524.6667
$ws.Range("N3").Value = -5127.6665
$ws.Range("H20").Value = 7748.4287
$ws.Range("I20").Value = 9997.5
$ws.Range("J20").Value = 4749.6665
$ws.Range("K20").Value = 9997.5
$ws.Range("L20").Value = 4749.6665
$ws.Range("M20").Value = -9750.5
$ws.Range("N20").Value = -5243.6665
$ws.Range("H37").Value = 750
$ws.Range("I37").Value = 750
$ws.Range("K37").Value = 750
$ws.Range("M37").Value = -613
$ws.Range("H86").Value = 9201.75
$ws.Range("I86").Value = 9000
$ws.Range("J86").Value = 9403.5
$ws.Range("K86").Value = 9000
$ws.Range("L86").Value = 9403.5
$ws.Range("M86").Value = -7877
$ws.Range("N86").Value = -11649.5
$ws.Range("H89").Value = 9201.75
$ws.Range("I89").Value = 9000
$ws.Range("J89").Value = 9403.5
$ws.Range("K89").Value = 45000
$ws.Range("L89").Value = 47017.5
$ws.Range("M89").Value = -39384
$ws.Range("N89").Value = -58249.5
$ws.Range("H94").Value = 922.2222
$ws.Range("I94").Value = 950.3077
$ws.Range("J94").Value = 849.2
$ws.Range("K94").Value = 950.3077
$ws.Range("L94").Value = 849.2
$ws.Range("M94").Value = -499.3077
$ws.Range("N94").Value = -1751.2
$ws.Range("H107").Value = 1164.8572
$ws.Range("I107").Value = 1159
$ws.Range("K107").Value = 1159
$ws.Range("M107").Value = 761
$ws.Range("H122").Value = 1979899
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H133").Value = 150000
$ws.Range("J133").Value = 150000
$ws.Range("L133").Value = 150000
$ws.Range("N133").Value = -160120

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 16899.8
$ws.Range("J28").Value = 18624.75
$ws.Range("L28").Value = 18624.75
$ws.Range("N28").Value = -19114.75
$ws.Range("H31").Value = 3687.7778
$ws.Range("J31").Value = 4163.6
$ws.Range("L31").Value = 4163.6
$ws.Range("N31").Value = -4753.6
$ws.Range("H34").Value = 3687.7778
$ws.Range("J34").Value = 4163.6
$ws.Range("L34").Value = 4163.6
$ws.Range("N34").Value = -4567.6
$ws.Range("H35").Value = 12497.5
$ws.Range("I35").Value = 12497.5
$ws.Range("K35").Value = 12497.5
$ws.Range("M35").Value = -12203.5
$ws.Range("H134").Value = 2497.45
$ws.Range("I134").Value = 2051.1072
$ws.Range("K134").Value = 6153.321599999999
$ws.Range("M134").Value = -3618.321599999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 898.5
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 898.5
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 2695.5
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -5315.5
$ws.Range("H132").Value = 10614.077
$ws.Range("I132").Value = 10543.909
$ws.Range("J132").Value = 11000
$ws.Range("K132").Value = 94895.181
$ws.Range("L132").Value = 99000
$ws.Range("M132").Value = -92365.181
$ws.Range("N132").Value = -104060

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4751.0835
$ws.Range("I80").Value = 4544.1665
$ws.Range("J80").Value = 4958
$ws.Range("K80").Value = 4544.1665
$ws.Range("L80").Value = 4958
$ws.Range("M80").Value = -3546.1665
$ws.Range("N80").Value = -6954
$ws.Range("H83").Value = 4751.0835
$ws.Range("I83").Value = 4544.1665
$ws.Range("J83").Value = 4958
$ws.Range("K83").Value = 22720.8325
$ws.Range("L83").Value = 24790
$ws.Range("M83").Value = -17728.8325
$ws.Range("N83").Value = -34774
$ws.Range("H128").Value = 89999
$ws.Range("J128").Value = 89999
$ws.Range("L128").Value = 89999
$ws.Range("N128").Value = -99959

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 4900
$ws.Range("I32").Value = 4900
$ws.Range("K32").Value = 4900
$ws.Range("M32").Value = -4583
$ws.Range("H40").Value = 3461.75
$ws.Range("I40").Value = 3527.7144
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 3527.7144
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -3391.7144
$ws.Range("N40").Value = -3272
$ws.Range("H46").Value = 2848.2173
$ws.Range("I46").Value = 2063.1875
$ws.Range("K46").Value = 2063.1875
$ws.Range("M46").Value = -1875.1875

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 13084.053
$ws.Range("J81").Value = 12162.4375
$ws.Range("L81").Value = 24324.875
$ws.Range("N81").Value = -26446.875
$ws.Range("H84").Value = 13084.053
$ws.Range("J84").Value = 12162.4375
$ws.Range("L84").Value = 121624.375
$ws.Range("N84").Value = -132232.375
$ws.Range("H99").Value = 105000
$ws.Range("I99").Value = 150000
$ws.Range("J99").Value = 82500
$ws.Range("K99").Value = 150000
$ws.Range("L99").Value = 82500
$ws.Range("M99").Value = -147005
$ws.Range("N99").Value = -88490
$ws.Range("H113").Value = 1263.5
$ws.Range("I113").Value = 1116.2
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 3348.6
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -1178.6
$ws.Range("N113").Value = -10340
$ws.Range("H136").Value = 1510.1111
$ws.Range("I136").Value = 1510.1111
$ws.Range("K136").Value = 4530.3333
$ws.Range("M136").Value = -1980.3333
